$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 126 - this shifts existing rows 126:154 down to 127:155
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are identical across all rows in this block
# (same market/product/category), mirroring row 125 and the other rows.
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C126").Value = "Ñuble"
$ws.Range("D126").Value = 44504
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = "Fruta"
$ws.Range("G126").Value = 100104
$ws.Range("H126").Value = "Frutos de pepita"
$ws.Range("I126").Value = 100104005
$ws.Range("J126").Value = "Pera"
$ws.Range("K126").Value = "Packham's Triumph"
$ws.Range("L126").Value = "Primera"
$ws.Range("M126").Value = 120
$ws.Range("N126").Value = 10000
$ws.Range("O126").Value = 11000
$ws.Range("P126").Value = 10500
$ws.Range("Q126").Value = "$/caja 16 kilos empedrada"
$ws.Range("R126").Value = "Provincia de Curicó"
$ws.Range("S126").Value = 656
$ws.Range("T126").Value = 16

# Match the date number format used by the other cells in column D (style index 2)
$ws.Range("D126").NumberFormat = $ws.Range("D125").NumberFormat
